# "Updated Core API Tests"
#
# The GNUKhata test-plan sheet drops its "Running Organisation Suite" test
# row (old row 12) entirely - deleting the whole row shifts every row below
# it up by one, which is why row 12's old neighbours inherit each other's
# formatting further down the sheet (rows 13-49) and the former last row
# (50) disappears. Separately, the G column ("story / feature" pass/fail
# flag) for the still-current suites (rows 5-11) flips from "yes" to "no".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content edit: G5:G11 "yes" -> "no" --------------------------------
$ws.Range("G5:G11").Value = "no"

# --- Structural edit: delete the "Running Organisation Suite" row ------
# (old row 12 - "Organisation" / "Organisation.data"). Everything below
# shifts up by one row, the sheet shrinks from 50 to 49 data rows, and the
# now-orphaned shared strings ("Running Organisation Suite", "Organisation",
# "Organisation.data") are dropped automatically on save.
$ws.Rows(12).Delete()

# conditionalFormatting / dataValidation ranges that referenced the old
# last row (...:M50 / ...:H50) need to be re-anchored to the new extent
# (...:M49 / ...:H49) - Excel auto-shrinks the data validation sqref on
# row delete but leaves the "no blanks" conditional format range stale,
# so fix that one up explicitly.
$fcs = $ws.Range("A1:M49").FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.AppliesTo.Address() -eq '$A$5:$M$50') {
        $fc.ModifyAppliesToRange($ws.Range("A5:M49"))
    }
}

# --- View state: scroll / selection -------------------------------------
# Previously the window was scrolled to B2 with C12 selected (the row that
# is now gone); re-home the view near the top-left and select the cell
# that took over the old selection's role.
$excel.ActiveWindow.ScrollRow = 2
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E11").Select()
